$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: "Relations de 2e niveau;Relations de 3e niveau et plus" -> "Relations de 3e niveau et plus"
$ws.Range("B2").Value = "Relations de 3e niveau et plus"

# E2: "Français" -> cleared (empty cell, keeps style)
$ws.Range("E2").ClearContents()

# J2: "Ressources humaines" -> " " (single space)
$ws.Range("J2").Value = " "

# Update the window view: scroll so column L is the left-most visible column,
# and change the active selection from O10 to E2.
$win = $excel.ActiveWindow
$win.ScrollColumn = 12
$win.ScrollRow = 1
$ws.Range("E2").Select()
